$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ------------------------------------------------------------------
# 1. Copy the formatting (date style in col A, grey-fill style in
#    B:D) of the last existing week (row 280) down across the seven
#    new weekly blocks (each block = 7 day-rows, separated by one
#    blank spacer row, exactly like the existing data above it).
# ------------------------------------------------------------------
$ws.Range("A280:D280").Copy() | Out-Null
$ws.Range("A282:D288").PasteSpecial(-4122)   | Out-Null
$ws.Range("A290:D296").PasteSpecial(-4122)   | Out-Null
$ws.Range("A298:D304").PasteSpecial(-4122)   | Out-Null
$ws.Range("A306:D312").PasteSpecial(-4122)   | Out-Null
$ws.Range("A314:D320").PasteSpecial(-4122)   | Out-Null
$ws.Range("A322:D328").PasteSpecial(-4122)   | Out-Null
$ws.Range("A330:D336").PasteSpecial(-4122)   | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. The old single-row "totals" block used to live at row 283
#    (right after the last populated day, row 280). That content
#    is being replaced by real calendar days now, so the stray
#    formulas in B283/C283/D283 need to go before row 283 becomes
#    an ordinary data row again.
# ------------------------------------------------------------------
$ws.Range("B283:D283").ClearContents() | Out-Null

# ------------------------------------------------------------------
# 3. Fill in column A (the running date, one per calendar day) for
#    every day row from 2023-05-08 (serial 45054) through
#    2023-06-25 (serial 45102), skipping the blank spacer rows that
#    separate each 7-day week (289, 297, 305, 313, 321, 329, 337 -
#    those stay empty, same as every earlier week boundary).
# ------------------------------------------------------------------
$weekStartRows = @(282, 290, 298, 306, 314, 322, 330)
$serial = 45054
foreach ($startRow in $weekStartRows) {
    for ($i = 0; $i -lt 7; $i++) {
        $r = $startRow + $i
        $ws.Cells.Item($r, 1).Value = $serial
        $serial = $serial + 1
    }
}

# ------------------------------------------------------------------
# 4. Only the first new week (already lived through) has real daily
#    counts recorded in column B; every following week is still
#    blank, waiting to be filled in day by day.
# ------------------------------------------------------------------
$firstWeekCounts = @(5, 7, 6, 4, 11, 6)
for ($i = 0; $i -lt $firstWeekCounts.Length; $i++) {
    $ws.Cells.Item(282 + $i, 2).Value = $firstWeekCounts[$i]
}

# ------------------------------------------------------------------
# 5. Re-create the running totals, now three rows further down (338)
#    with one extra blank spacer row (337) above them, exactly like
#    every other week boundary.
# ------------------------------------------------------------------
$ws.Range("B280").Copy() | Out-Null
$ws.Range("B338").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B338").Formula = "=SUM(B2:B287)"
$ws.Range("C338").Formula = "=B338/31"
$ws.Range("D338").Formula = "=35*7"

# ------------------------------------------------------------------
# 6. Restore the view state: selection on D286 (as in the saved
#    file) with the window scrolled near the bottom of the newly
#    added data.
# ------------------------------------------------------------------
$ws.Range("D286").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 313
